$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed symbol list (Price / Volume(1h)% / Hora) as scraped 2023-01-16.
# These columns are stored as literal text in the workbook (matching the
# original inline-string cells), so each value is entered with a leading
# apostrophe to force text entry and stop Excel from re-parsing numbers /
# percentages. The style is reset back to Normal afterwards so the text
# entry doesn't leave a stray number-format behind on the cell.
$updates = @(
    @{ Cell = "D2"; Value = "304.77" },
    @{ Cell = "E2"; Value = "2.08%" },
    @{ Cell = "G2"; Value = "3" },
    @{ Cell = "D3"; Value = "32.08" },
    @{ Cell = "E3"; Value = "1.88%" },
    @{ Cell = "G3"; Value = "3" },
    @{ Cell = "D4"; Value = "5.194" },
    @{ Cell = "E4"; Value = "1.04%" },
    @{ Cell = "G4"; Value = "3" },
    @{ Cell = "D5"; Value = "0.07459" },
    @{ Cell = "E5"; Value = "-0.45%" },
    @{ Cell = "G5"; Value = "3" },
    @{ Cell = "D6"; Value = "2.379" },
    @{ Cell = "E6"; Value = "44.47%" },
    @{ Cell = "G6"; Value = "3" },
    @{ Cell = "D7"; Value = "7.990" },
    @{ Cell = "E7"; Value = "2.08%" },
    @{ Cell = "G7"; Value = "3" },
    @{ Cell = "D8"; Value = "3.867" },
    @{ Cell = "E8"; Value = "1.30%" },
    @{ Cell = "G8"; Value = "3" },
    @{ Cell = "D9"; Value = "0.9203" },
    @{ Cell = "E9"; Value = "-0.29%" },
    @{ Cell = "G9"; Value = "3" },
    @{ Cell = "D10"; Value = "0.1730" },
    @{ Cell = "E10"; Value = "0.88%" },
    @{ Cell = "G10"; Value = "3" },
    @{ Cell = "D11"; Value = "0.07661" },
    @{ Cell = "E11"; Value = "-0.25%" },
    @{ Cell = "G11"; Value = "3" },
    @{ Cell = "D12"; Value = "0.08212" },
    @{ Cell = "E12"; Value = "2.53%" },
    @{ Cell = "G12"; Value = "3" },
    @{ Cell = "D13"; Value = "0.03010" },
    @{ Cell = "E13"; Value = "0.32%" },
    @{ Cell = "G13"; Value = "3" },
    @{ Cell = "D14"; Value = "0.09933" },
    @{ Cell = "E14"; Value = "0.25%" },
    @{ Cell = "G14"; Value = "3" },
    @{ Cell = "D15"; Value = "0.001511" },
    @{ Cell = "E15"; Value = "0.64%" },
    @{ Cell = "G15"; Value = "3" },
    @{ Cell = "D16"; Value = "0.006116" },
    @{ Cell = "E16"; Value = "-3.25%" },
    @{ Cell = "G16"; Value = "3" },
    @{ Cell = "D17"; Value = "3.495" },
    @{ Cell = "E17"; Value = "1.44%" },
    @{ Cell = "G17"; Value = "3" },
    @{ Cell = "D18"; Value = "2.228" },
    @{ Cell = "E18"; Value = "0.03%" },
    @{ Cell = "G18"; Value = "3" },
    @{ Cell = "D19"; Value = "0.3261" },
    @{ Cell = "E19"; Value = "-0.97%" },
    @{ Cell = "G19"; Value = "3" },
    @{ Cell = "D20"; Value = "0.1348" },
    @{ Cell = "E20"; Value = "-0.11%" },
    @{ Cell = "G20"; Value = "3" },
    @{ Cell = "D21"; Value = "4.660" },
    @{ Cell = "E21"; Value = "2.04%" },
    @{ Cell = "G21"; Value = "3" },
    @{ Cell = "D22"; Value = "0.04607" },
    @{ Cell = "E22"; Value = "-1.43%" },
    @{ Cell = "G22"; Value = "3" },
    @{ Cell = "D23"; Value = "0.1565" },
    @{ Cell = "E23"; Value = "0.79%" },
    @{ Cell = "G23"; Value = "3" },
    @{ Cell = "D24"; Value = "0.001259" },
    @{ Cell = "E24"; Value = "3.06%" },
    @{ Cell = "G24"; Value = "3" },
    @{ Cell = "D25"; Value = "0.004519" },
    @{ Cell = "E25"; Value = "2.16%" },
    @{ Cell = "G25"; Value = "3" },
    @{ Cell = "D26"; Value = "0.0001299" },
    @{ Cell = "E26"; Value = "-7.30%" },
    @{ Cell = "G26"; Value = "3" },
    @{ Cell = "D27"; Value = "0.0002740" },
    @{ Cell = "E27"; Value = "52.12%" },
    @{ Cell = "G27"; Value = "3" },
    @{ Cell = "G28"; Value = "3" },
    @{ Cell = "G29"; Value = "3" },
    @{ Cell = "G30"; Value = "3" },
    @{ Cell = "G31"; Value = "3" },
    @{ Cell = "G32"; Value = "3" },
    @{ Cell = "G33"; Value = "3" },
    @{ Cell = "G34"; Value = "3" },
    @{ Cell = "G35"; Value = "3" },
    @{ Cell = "G36"; Value = "3" },
    @{ Cell = "G37"; Value = "3" },
    @{ Cell = "G38"; Value = "3" },
    @{ Cell = "D39"; Value = "0.01774" },
    @{ Cell = "E39"; Value = "7.55%" },
    @{ Cell = "G39"; Value = "3" },
    @{ Cell = "D40"; Value = "0.04553" },
    @{ Cell = "E40"; Value = "0.08%" },
    @{ Cell = "G40"; Value = "3" },
    @{ Cell = "D41"; Value = "0.007372" },
    @{ Cell = "E41"; Value = "6.24%" },
    @{ Cell = "G41"; Value = "3" },
    @{ Cell = "D42"; Value = "0.1362" },
    @{ Cell = "E42"; Value = "1.35%" },
    @{ Cell = "G42"; Value = "3" },
    @{ Cell = "D43"; Value = "0.002178" },
    @{ Cell = "E43"; Value = "5.65%" },
    @{ Cell = "G43"; Value = "3" },
    @{ Cell = "E44"; Value = "-19.30%" },
    @{ Cell = "G44"; Value = "3" },
    @{ Cell = "D45"; Value = "0.00006110" },
    @{ Cell = "E45"; Value = "0.04%" },
    @{ Cell = "G45"; Value = "3" },
    @{ Cell = "E46"; Value = "-57.20%" },
    @{ Cell = "G46"; Value = "3" },
    @{ Cell = "D47"; Value = "0.009892" },
    @{ Cell = "E47"; Value = "-19.42%" },
    @{ Cell = "G47"; Value = "3" },
    @{ Cell = "G48"; Value = "3" },
    @{ Cell = "G49"; Value = "3" },
    @{ Cell = "G50"; Value = "3" },
    @{ Cell = "G51"; Value = "3" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.Value = "'" + $u.Value
    $cell.Style = "Normal"
}

